{"js": "// Replace the 25 \"A\u00f7B=C, D\" division-answer cells with their new values.\n// Each old value is unique in the document, so an exact (non-wildcard)\n// search-and-replace per pair is sufficient and keeps run formatting\n// (font / size) untouched because insertText(\"Replace\") reuses the\n// existing run's formatting.\nconst replacements = [\n  [\"818\u00f72=409, 0\", \"740\u00f77=105, 5\"],\n  [\"288\u00f76=48, 0\", \"645\u00f75=129, 0\"],\n  [\"828\u00f75=165, 3\", \"111\u00f74=27, 3\"],\n  [\"541\u00f78=67, 5\", \"979\u00f79=108, 7\"],\n  [\"230\u00f74=57, 2\", \"202\u00f72=101, 0\"],\n  [\"557\u00f76=92, 5\", \"502\u00f78=62, 6\"],\n  [\"884\u00f77=126, 2\", \"645\u00f73=215, 0\"],\n  [\"375\u00f73=125, 0\", \"670\u00f75=134, 0\"],\n  [\"891\u00f77=127, 2\", \"557\u00f78=69, 5\"],\n  [\"137\u00f77=19, 4\", \"882\u00f73=294, 0\"],\n  [\"793\u00f79=88, 1\", \"442\u00f76=73, 4\"],\n  [\"714\u00f75=142, 4\", \"501\u00f74=125, 1\"],\n  [\"674\u00f79=74, 8\", \"849\u00f74=212, 1\"],\n  [\"161\u00f72=80, 1\", \"738\u00f73=246, 0\"],\n  [\"704\u00f72=352, 0\", \"292\u00f77=41, 5\"],\n  [\"682\u00f79=75, 7\", \"167\u00f78=20, 7\"],\n  [\"185\u00f74=46, 1\", \"952\u00f75=190, 2\"],\n  [\"797\u00f74=199, 1\", \"826\u00f74=206, 2\"],\n  [\"318\u00f78=39, 6\", \"683\u00f73=227, 2\"],\n  [\"706\u00f79=78, 4\", \"406\u00f79=45, 1\"],\n  [\"702\u00f75=140, 2\", \"132\u00f77=18, 6\"],\n  [\"351\u00f76=58, 3\", \"187\u00f72=93, 1\"],\n  [\"784\u00f74=196, 0\", \"951\u00f72=475, 1\"],\n  [\"469\u00f79=52, 1\", \"753\u00f73=251, 0\"],\n  [\"562\u00f72=281, 0\", \"763\u00f78=95, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"A\u00f7B=C, D\" division-answer cells with their new values.\n# Every old value is unique in the document, so Find/Replace (no wildcards,\n# match-case) of the whole cell text is sufficient; it reuses the existing\n# run's formatting (font / size) so nothing else in the document changes.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"818\u00f72=409, 0\", \"740\u00f77=105, 5\"),\n    @(\"288\u00f76=48, 0\", \"645\u00f75=129, 0\"),\n    @(\"828\u00f75=165, 3\", \"111\u00f74=27, 3\"),\n    @(\"541\u00f78=67, 5\", \"979\u00f79=108, 7\"),\n    @(\"230\u00f74=57, 2\", \"202\u00f72=101, 0\"),\n    @(\"557\u00f76=92, 5\", \"502\u00f78=62, 6\"),\n    @(\"884\u00f77=126, 2\", \"645\u00f73=215, 0\"),\n    @(\"375\u00f73=125, 0\", \"670\u00f75=134, 0\"),\n    @(\"891\u00f77=127, 2\", \"557\u00f78=69, 5\"),\n    @(\"137\u00f77=19, 4\", \"882\u00f73=294, 0\"),\n    @(\"793\u00f79=88, 1\", \"442\u00f76=73, 4\"),\n    @(\"714\u00f75=142, 4\", \"501\u00f74=125, 1\"),\n    @(\"674\u00f79=74, 8\", \"849\u00f74=212, 1\"),\n    @(\"161\u00f72=80, 1\", \"738\u00f73=246, 0\"),\n    @(\"704\u00f72=352, 0\", \"292\u00f77=41, 5\"),\n    @(\"682\u00f79=75, 7\", \"167\u00f78=20, 7\"),\n    @(\"185\u00f74=46, 1\", \"952\u00f75=190, 2\"),\n    @(\"797\u00f74=199, 1\", \"826\u00f74=206, 2\"),\n    @(\"318\u00f78=39, 6\", \"683\u00f73=227, 2\"),\n    @(\"706\u00f79=78, 4\", \"406\u00f79=45, 1\"),\n    @(\"702\u00f75=140, 2\", \"132\u00f77=18, 6\"),\n    @(\"351\u00f76=58, 3\", \"187\u00f72=93, 1\"),\n    @(\"784\u00f74=196, 0\", \"951\u00f72=475, 1\"),\n    @(\"469\u00f79=52, 1\", \"753\u00f73=251, 0\"),\n    @(\"562\u00f72=281, 0\", \"763\u00f78=95, 3\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
